$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9:19 down to 10:20
$ws.Rows.Item(9).Insert()

# Set the new row 9 values
$ws.Cells.Item(9, 1).Value = 50
$ws.Cells.Item(9, 2).Value = "no water"

# Update row 10 (previously row 9) A value from 10 to 55
$ws.Cells.Item(10, 1).Value = 55

# Update selection to B9
$ws.Range("B9").Select()
